# Fruta / hortaliza, semanal
# Insert a new weekly record as the first data row (row 6), pushing the
# existing rows 6-12 down to rows 7-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6..12 down to 7..13, keeping row 6 free (and its date style)
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new weekly reading
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value2 = 44512
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112022
$ws.Range("G6").Value = "Arveja Verde"
$ws.Range("H6").Value = "Perfection"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 14500
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 580
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
